# Auto-generated Excel COM-interop edit script
# Updates "想去人数" (want-to-go count) figures and one event cover image URL
# across the "展览" (sheet 1), "演出" (sheet 2) and "全部类型" (sheet 4) worksheets,
# matching a refreshed data scrape (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Worksheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1233  # was 1230
$ws.Range("F3").Value = 1124  # was 1122
$ws.Range("F4").Value = 874  # was 872
$ws.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202402/uPDIsIoV1708311822716.jpeg"  # was "//i1.hdslb.com/bfs/openplatform/202402/IWtNSAXt1707014699653.jpeg"
$ws.Range("F7").Value = 634  # was 632
$ws.Range("F8").Value = 83  # was 82
$ws.Range("F9").Value = 42  # was 41
$ws.Range("F11").Value = 2235  # was 2223
$ws.Range("F12").Value = 1555  # was 1552
$ws.Range("F13").Value = 1246  # was 1240
$ws.Range("F17").Value = 720  # was 715
$ws.Range("F18").Value = 266  # was 265
$ws.Range("F19").Value = 1082  # was 1080
$ws.Range("F22").Value = 4222  # was 4194
$ws.Range("F23").Value = 205  # was 204
$ws.Range("F24").Value = 144  # was 143
$ws.Range("F28").Value = 603  # was 602
$ws.Range("F33").Value = 356  # was 355
$ws.Range("F34").Value = 922  # was 921
$ws.Range("F36").Value = 87  # was 86
$ws.Range("F37").Value = 113  # was 112

# --- Worksheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 15  # was 14

# --- Worksheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1234  # was 1230
$ws.Range("F5").Value = 1124  # was 1122
$ws.Range("F6").Value = 874  # was 872
$ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202402/uPDIsIoV1708311822716.jpeg"  # was "//i1.hdslb.com/bfs/openplatform/202402/IWtNSAXt1707014699653.jpeg"
$ws.Range("F11").Value = 634  # was 632
$ws.Range("F12").Value = 83  # was 82
$ws.Range("F13").Value = 42  # was 41
$ws.Range("F16").Value = 2235  # was 2223
$ws.Range("F17").Value = 1555  # was 1552
$ws.Range("F18").Value = 1246  # was 1240
$ws.Range("F23").Value = 720  # was 715
$ws.Range("F24").Value = 266  # was 265
$ws.Range("F25").Value = 1082  # was 1080
$ws.Range("F28").Value = 4222  # was 4194
$ws.Range("F29").Value = 205  # was 204
$ws.Range("F30").Value = 144  # was 143
$ws.Range("F34").Value = 603  # was 602
$ws.Range("F39").Value = 356  # was 355
$ws.Range("F40").Value = 922  # was 921
$ws.Range("F42").Value = 87  # was 86
$ws.Range("F43").Value = 113  # was 112
$ws.Range("F46").Value = 15  # was 14
